$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 277
$ws.Range("B3").Value = 719
$ws.Range("B5").Value = 635
$ws.Range("B9").Value = 510
$ws.Range("B10").Value = 496
$ws.Range("B11").Value = 581
$ws.Range("B15").Value = 456
$ws.Range("B16").Value = 424
$ws.Range("B18").Value = 455
$ws.Range("B19").Value = 389
$ws.Range("B41").Value = 480
$ws.Range("B42").Value = 272

$ws.Range("A43").Value = 42
$ws.Range("B43").Value = 3
